# Cosmetic Changes.xlsx — add a new "extra comment" column (F) on the
# "Make objects kickable" row and bump its Priority from "Low" to
# "Low (very high)".
#
# NOTE: the original author's commit mostly touched raw OOXML plumbing
# (fileVersion build number, absPath, revisionPtr GUID, window geometry,
# calcFeatures list) that simply reflects which machine/Excel build saved
# the file - none of that is reachable (or meaningful) through the Excel
# object model, so this script focuses purely on the actual worksheet
# content edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Order matters: typing into F10 first, then D10, reproduces the shared
# string table order of the target workbook (new strings are appended in
# the order they're entered).
$ws.Range("F10").Value2 = "LETS GOOOOOOO FOOTBALL MINIGAME WITH MUSHROOMS YEAHHH"
$ws.Range("D10").Value2 = "Low (very high)"

# Give the new comment cell a thin right-hand border (matches the rest of
# the table's boxed look) without touching its neighbours.
$border = $ws.Range("F10").Borders.Item(10)
$border.LineStyle = 1
$border.Weight = 2
